$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-02-15 Saturday" "2025-02-16 Sunday"
Replace-Text "454÷7=" "866÷5="
Replace-Text "108÷3=" "663÷3="
Replace-Text "858÷2=" "868÷5="
Replace-Text "394÷9=" "741÷2="
Replace-Text "975÷2=" "993÷5="
Replace-Text "598÷3=" "714÷6="
Replace-Text "544÷9=" "205÷7="
Replace-Text "395÷9=" "301÷8="
Replace-Text "919÷6=" "547÷6="
Replace-Text "763÷3=" "763÷4="
Replace-Text "188÷3=" "656÷6="
Replace-Text "219÷3=" "745÷7="
Replace-Text "990÷9=" "126÷7="
Replace-Text "523÷3=" "230÷3="
Replace-Text "631÷7=" "743÷5="
Replace-Text "981÷3=" "168÷2="
Replace-Text "455÷9=" "177÷4="
Replace-Text "311÷6=" "388÷7="
Replace-Text "555÷9=" "831÷9="
Replace-Text "114÷5=" "523÷5="
Replace-Text "775÷7=" "291÷6="
Replace-Text "475÷3=" "341÷7="
Replace-Text "857÷6=" "406÷9="
Replace-Text "322÷8=" "394÷5="
Replace-Text "980÷6=" "552÷9="
